# Fixes to launch the browser
# Adds two new test-case sheets (facebook_Login_TC002, Demo_TC_TC003) and
# registers them as two new rows on the Index sheet.

$wb = $excel.ActiveWorkbook

$wsIndex = $wb.Worksheets.Item("Index")
$wsTC001 = $wb.Worksheets.Item("facebook_Login_TC001")

# ---------------------------------------------------------------------
# 1. Index sheet - append rows for the two new test cases
# ---------------------------------------------------------------------

# Clone the formatting of the existing data row (row 2) down into rows 3-4
$wsIndex.Range("A2:F2").Copy()
$wsIndex.Range("A3:F4").PasteSpecial(-4122)   # xlPasteFormats
$wsIndex.Rows("3").RowHeight = $wsIndex.Rows("2").RowHeight
$wsIndex.Rows("4").RowHeight = $wsIndex.Rows("2").RowHeight

$wsIndex.Range("A3").Value = "UI_Validation_Testcases"
$wsIndex.Range("B3").Value = "facebook_Login_TC002"
$wsIndex.Range("C3").Value = "Yes"
$wsIndex.Range("D3").Value = "Prod"
$wsIndex.Range("E3").Value = "Verify facebook SignUp page"
$wsIndex.Range("F3").Value = "Smoke"

$wsIndex.Range("A4").Value = "UI_Validation_Testcases"
$wsIndex.Range("B4").Value = "Demo_TC_TC003"
$wsIndex.Range("C4").Value = "Yes"
$wsIndex.Range("D4").Value = "Prod"
$wsIndex.Range("E4").Value = "What done is done when we say its done"
$wsIndex.Range("F4").Value = "Dummy"

$wsIndex.Columns("E").ColumnWidth = 43.13

$wsIndex.Range("C3:C4").Select()

# ---------------------------------------------------------------------
# 2. facebook_Login_TC002 sheet - copy of facebook_Login_TC001 plus a
#    FirstName column
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTC002 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsTC002.Name = "facebook_Login_TC002"

$wsTC001.Range("A1:B2").Copy()
$wsTC002.Range("A1:B2").PasteSpecial(-4122)   # xlPasteFormats

$wsTC002.Columns("A").ColumnWidth = 18.72
$wsTC002.Columns("B").ColumnWidth = 9.92
$wsTC002.Columns("C").ColumnWidth = 9.71

$wsTC002.Range("A1").Value = "Username"
$wsTC002.Range("B1").Value = "Password"
$wsTC002.Range("A2").Value = "tonystark@gmail.com"
$wsTC002.Range("B2").Value = "tony56432"

$wsTC002.Range("B1").Copy()
$wsTC002.Range("C1").PasteSpecial(-4122)
$wsTC002.Range("C1").Value = "FirstName"

$wsTC002.Range("B2").Copy()
$wsTC002.Range("C2").PasteSpecial(-4122)
$wsTC002.Range("C2").Value = "David"

$wsTC002.Hyperlinks.Add($wsTC002.Range("A2"), "mailto:tonystark@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "tonystark@gmail.com")

# Re-apply the TC001 hyperlink-cell formatting (font/style), which the
# Hyperlinks.Add call above would otherwise overwrite with a generic
# "Hyperlink" style.
$wsTC001.Range("A2").Copy()
$wsTC002.Range("A2").PasteSpecial(-4122)
$wsTC002.Range("A2").Value = "tonystark@gmail.com"

$wsTC002.Range("C1").Select()

# ---------------------------------------------------------------------
# 3. Demo_TC_TC003 sheet - placeholder test case
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTC003 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsTC003.Name = "Demo_TC_TC003"

$wsTC003.Range("A1").Value = "No Data"
$wsTC003.Range("A2").Value = "No Data"

$wsTC003.Range("A1").Select()

$wsIndex.Activate()
